$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Insert a new "comments" column right after column A, shifting the
# existing wave columns (old B..AF) one to the right (new C..AG).
$ws.Columns("B:B").Insert()

# The old last data row (A31 "fehlender Burgernahe der Politiker", with
# the same two values already present on the "Uberforderung der
# Politiker" row) is now redundant -- drop it.
$ws.Rows("31:31").Delete()

# Give the new comments column a wider, custom width (close match to the
# author's manual resize; the COM width setter here only supports ~1/6
# character-unit granularity).
$ws.Columns("B:B").ColumnWidth = 15.65

# Header for the new column.
$ws.Range("B1").Value = "comments"

# Row-specific notes explaining wording/coverage changes over the waves.
$ws.Range("B14").Value = "item wording changed from ""fehlender Bürgernähe der Politiker"" to ""Überforderung der Politiker"" possibly from w16 onwards."
$ws.Range("B16").Value = "item wording changed to included Wetterextreme in later waves, unclear which wave exactly. Possible overlap with Klimawandel."
$ws.Range("B2").Value = "inclusion of Corona-Infektion from w29 onwards."

# Make "data" the active/selected sheet (was "dates"), with A11 selected.
$ws.Range("A11").Select()
$ws.Activate()
